# Wydzielenie kodu do osobnych klas
#
# Refreshes the "exclusion window" date ranges on a handful of sheets
# (shifting old 2018 sample dates forward to 2024), adds a new blank
# row to WYKLUCZENIA_MASZYNY, and updates the in-workbook selection /
# active-sheet UI state to match where the author left off editing.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# WYKLUCZENIA_PRACOWNICY - refresh exclusion date ranges
# ---------------------------------------------------------------
$wsPrac = $wb.Worksheets.Item("WYKLUCZENIA_PRACOWNICY")
$wsPrac.Range("B2").Value = 45499
$wsPrac.Range("C2").Value = 45504
$wsPrac.Range("B3").Value = 45441
$wsPrac.Range("C3").Value = 45443
$wsPrac.Range("B4").Value = 45536
$wsPrac.Range("C4").Value = 45539
$wsPrac.Range("B5").Value = 45542
$wsPrac.Range("C5").Value = 45546
$wsPrac.Range("B6").Value = 45468
$wsPrac.Range("C6").Value = 45470

# ---------------------------------------------------------------
# ZAMOWIENIA_KLIENTA - refresh a couple of order due-dates
# ---------------------------------------------------------------
$wsZam = $wb.Worksheets.Item("ZAMOWIENIA_KLIENTA")
$wsZam.Range("C4").Value = 45639
$wsZam.Range("C10").Value = 45470

# ---------------------------------------------------------------
# WYKLUCZENIA_MASZYNY - refresh dates and append a new blank row
# ---------------------------------------------------------------
$wsMasz = $wb.Worksheets.Item("WYKLUCZENIA_MASZYNY")
$wsMasz.Range("B2").Value = 45466
$wsMasz.Range("C2").Value = 45471
$wsMasz.Range("B3").Value = 45573
$wsMasz.Range("C3").Value = 45575
$wsMasz.Range("B4").Value = 45631
$wsMasz.Range("C4").Value = 45639
$wsMasz.Range("B5").Value = 45605
$wsMasz.Range("C5").Value = 45606
$wsMasz.Range("B6").Value = 45503
$wsMasz.Range("C6").Value = 45504

# New row 7: an empty, date-formatted cell in C7 (same number format
# as the rest of the DATA_DO column) - mirrors the author starting a
# new exclusion entry.
$wsMasz.Range("C6").Copy($wsMasz.Range("C7"))
$wsMasz.Range("C7").ClearContents()

# ---------------------------------------------------------------
# Leave the UI selection state the way the author left it: a cell
# selected on WYKLUCZENIA_PRACOWNICY, a range selected on
# ZAMOWIENIA_KLIENTA, and WYKLUCZENIA_MASZYNY as the final active
# sheet/selection (it must be the LAST sheet activated so it ends up
# as the saved active tab).
# ---------------------------------------------------------------
$wsPrac.Range("C8").Select()

$wsZam.Range("C6:C10").Select()

$wsMasz.Activate()
$wsMasz.Range("B5").Select()
